# Weekly price update: insert a new record row for "Albahaca" (Vega Modelo de
# Temuco) just before the existing row 249, shifting all subsequent rows down
# by one (dimension grows from A1:R308 to A1:R309).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 249 (native Excel Rows.Insert shifts
# rows 249..308 down to 250..309 and extends the used range accordingly).
$ws.Rows.Item(249).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(249, 1).Value  = 10
$ws.Cells.Item(249, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(249, 3).Value  = "La Araucanía"
$ws.Cells.Item(249, 4).Value  = 44932
$ws.Cells.Item(249, 5).Value  = 9
$ws.Cells.Item(249, 6).Value  = 100112052
$ws.Cells.Item(249, 7).Value  = "Albahaca"
$ws.Cells.Item(249, 8).Value  = "Sin especificar"
$ws.Cells.Item(249, 9).Value  = "Primera"
$ws.Cells.Item(249, 10).Value = 80
$ws.Cells.Item(249, 11).Value = 6000
$ws.Cells.Item(249, 12).Value = 6000
$ws.Cells.Item(249, 13).Value = 6000
$ws.Cells.Item(249, 14).Value = "`$/paquete"
$ws.Cells.Item(249, 15).Value = "Región del Maule"
$ws.Cells.Item(249, 16).Value = 6000
$ws.Cells.Item(249, 17).Value = 1
$ws.Cells.Item(249, 18).Value = "Hortaliza"
